# Daily auto-update: for each data row (2..last) in column E ("剩余" / days
# remaining) and F ("开始时间" / cycle-start date, stored as an integer
# YYYYMMDD), decrement the remaining-day counter by one. When the counter
# is already at 1 (last day of the cycle), instead roll it over to a fresh
# cycle: reset E back to the row's total-days value (column D) and advance
# F by that many days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # D: 总天 (total days)
    $eCell = $ws.Cells.Item($r, 5)   # E: 剩余 (days remaining)
    $fCell = $ws.Cells.Item($r, 6)   # F: 开始时间 (cycle-start date, yyyyMMdd)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }
    if ($dVal -eq "" -or $eVal -eq "" -or $fVal -eq "") {
        continue
    }

    $d = [int]$dVal
    $e = [int]$eVal
    $f = [int]$fVal

    if ($e -le 1) {
        # Cycle finished -> start a new one.
        $newE = $d

        $y = [int]([math]::Floor($f / 10000))
        $m = [int]([math]::Floor(($f % 10000) / 100))
        $day = [int]($f % 100)

        $dt = Get-Date -Year $y -Month $m -Day $day
        $dt = $dt.AddDays($d)

        $newF = [int]($dt.ToString("yyyyMMdd"))
    } else {
        $newE = $e - 1
        $newF = $f
    }

    $eCell.Value2 = $newE
    $fCell.Value2 = $newF
}
